# Event Expense - 23 July
# Remove the "Clone" action from the Approver Actions sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actions")

# Row 6 currently holds "Clone" - delete the whole row and shift the rows
# below it up (this removes the "Clone" entry from the Actions list while
# keeping "Delete" as the final entry).
$ws.Rows(6).Delete()

# Update the selection to match the new last row.
$ws.Range("A6:XFD6").Select()
